$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 394.1111
$ws.Range("I9").Value = 249.5
$ws.Range("K9").Value = 249.5
$ws.Range("M9").Value = -80.5
$ws.Range("H18").Value = 1616.5
$ws.Range("I18").Value = 1616.5
$ws.Range("K18").Value = 1616.5
$ws.Range("M18").Value = -1332.5
$ws.Range("H40").Value = 7252.4736
$ws.Range("I40").Value = 3945.1
$ws.Range("K40").Value = 3945.1
$ws.Range("M40").Value = -3770.1
$ws.Range("H42").Value = 150.25
$ws.Range("J42").Value = 282
$ws.Range("L42").Value = 846
$ws.Range("N42").Value = -1306
$ws.Range("H43").Value = 3149.5
$ws.Range("I43").Value = 3049.5
$ws.Range("K43").Value = 3049.5
$ws.Range("M43").Value = -2980.5
$ws.Range("H116").Value = 4173285
$ws.Range("I116").Value = 6537.5
$ws.Range("K116").Value = 6537.5
$ws.Range("M116").Value = -3095.5
$ws.Range("H135").Value = 1206.591
$ws.Range("I135").Value = 886.8421
$ws.Range("K135").Value = 7981.5789
$ws.Range("M135").Value = -5446.5789
$ws.Range("H137").Value = 331608.44
$ws.Range("I137").Value = 1662.0968
$ws.Range("J137").Value = 1118403.5
$ws.Range("K137").Value = 4986.2904
$ws.Range("L137").Value = 3355210.5
$ws.Range("M137").Value = -2436.2904
$ws.Range("N137").Value = -3360310.5
$ws.Range("H138").Value = 1741.7037
$ws.Range("J138").Value = 2874.25
$ws.Range("L138").Value = 8622.75
$ws.Range("N138").Value = -18902.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4938.279
$ws.Range("I32").Value = 1569.2885
$ws.Range("K32").Value = 1569.2885
$ws.Range("M32").Value = -1282.2885
$ws.Range("H39").Value = 21249.75
$ws.Range("I39").Value = 12500
$ws.Range("J39").Value = 29999.5
$ws.Range("K39").Value = 12500
$ws.Range("L39").Value = 29999.5
$ws.Range("M39").Value = -11980
$ws.Range("N39").Value = -31039.5
$ws.Range("H132").Value = 1808.5238
$ws.Range("I132").Value = 1634.027
$ws.Range("K132").Value = 4902.081
$ws.Range("M132").Value = -2372.081
$ws.Range("H138").Value = 30214.5
$ws.Range("J138").Value = 30214.5
$ws.Range("L138").Value = 30214.5
$ws.Range("N138").Value = -40494.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 80738
$ws.Range("I22").Value = 115124.78
$ws.Range("K22").Value = 115124.78
$ws.Range("M22").Value = -114951.78
$ws.Range("H99").Value = 3227580
$ws.Range("I99").Value = 335266.66
$ws.Range("J99").Value = 4467143
$ws.Range("K99").Value = 335266.66
$ws.Range("L99").Value = 4467143
$ws.Range("M99").Value = -333768.66
$ws.Range("N99").Value = -4470139
$ws.Range("H132").Value = 33620.277
$ws.Range("J132").Value = 33620.277
$ws.Range("L132").Value = 33620.277
$ws.Range("N132").Value = -43740.277
$ws.Range("H134").Value = 5407.8823
$ws.Range("I134").Value = 2316.9092
$ws.Range("J134").Value = 11074.667
$ws.Range("K134").Value = 6950.7276
$ws.Range("L134").Value = 33224.001
$ws.Range("M134").Value = -4415.7276
$ws.Range("N134").Value = -38294.001
$ws.Range("H140").Value = 43499
$ws.Range("J140").Value = 43499
$ws.Range("L140").Value = 43499
$ws.Range("N140").Value = -53859

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 683
$ws.Range("I22").Value = 558.6667
$ws.Range("J22").Value = 1149.25
$ws.Range("K22").Value = 558.6667
$ws.Range("L22").Value = 1149.25
$ws.Range("M22").Value = -208.6667
$ws.Range("N22").Value = -1849.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 900.5294
$ws.Range("I5").Value = 856.5
$ws.Range("J5").Value = 963.4286
$ws.Range("K5").Value = 2569.5
$ws.Range("L5").Value = 2890.2858
$ws.Range("M5").Value = -2457.5
$ws.Range("N5").Value = -3114.2858
$ws.Range("H113").Value = 2701993
$ws.Range("I113").Value = 1141.3334
$ws.Range("J113").Value = 4052419
$ws.Range("K113").Value = 3424.0002
$ws.Range("L113").Value = 12157257
$ws.Range("M113").Value = -1254.0002
$ws.Range("N113").Value = -12161597
$ws.Range("H122").Value = 777460.3
$ws.Range("J122").Value = 1443383.9
$ws.Range("L122").Value = 12990455.1
$ws.Range("N122").Value = -12995355.1
$ws.Range("H129").Value = 800
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H135").Value = 900.5294
$ws.Range("I135").Value = 856.5
$ws.Range("J135").Value = 963.4286
$ws.Range("K135").Value = 7708.5
$ws.Range("L135").Value = 8670.857399999999
$ws.Range("M135").Value = -5173.5
$ws.Range("N135").Value = -13740.8574

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2383517.5
$ws.Range("I113").Value = 2045.6666
$ws.Range("J113").Value = 6670167
$ws.Range("K113").Value = 2045.6666
$ws.Range("L113").Value = 6670167
$ws.Range("M113").Value = 124.3334
$ws.Range("N113").Value = -6674507
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()
$ws.Range("H132").Value = 3732.0833
$ws.Range("I132").Value = 3730.4546
$ws.Range("J132").Value = 3750
$ws.Range("K132").Value = 11191.3638
$ws.Range("L132").Value = 11250
$ws.Range("M132").Value = -8661.363799999999
$ws.Range("N132").Value = -16310

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H61").Value = 1399
$ws.Range("I61").Value = 998.75
$ws.Range("K61").Value = 998.75
$ws.Range("M61").Value = -796.75
$ws.Range("H113").Value = 1399
$ws.Range("I113").Value = 998.75
$ws.Range("K113").Value = 998.75
$ws.Range("M113").Value = 1171.25
$ws.Range("H136").Value = 1682.871
$ws.Range("J136").Value = 2612.3333
$ws.Range("L136").Value = 7836.999899999999
$ws.Range("N136").Value = -12936.9999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4999.143
$ws.Range("J62").Value = 4333
$ws.Range("L62").Value = 4333
$ws.Range("N62").Value = -5581
$ws.Range("H65").Value = 4999.143
$ws.Range("J65").Value = 4333
$ws.Range("L65").Value = 21665
$ws.Range("N65").Value = -27905
$ws.Range("H81").Value = 1625.5333
$ws.Range("I81").Value = 1037.2307
$ws.Range("J81").Value = 5449.5
$ws.Range("K81").Value = 2074.4614
$ws.Range("L81").Value = 10899
$ws.Range("M81").Value = -1013.4614
$ws.Range("N81").Value = -13021
$ws.Range("H84").Value = 1625.5333
$ws.Range("I84").Value = 1037.2307
$ws.Range("J84").Value = 5449.5
$ws.Range("K84").Value = 10372.307
$ws.Range("L84").Value = 54495
$ws.Range("M84").Value = -5068.307000000001
$ws.Range("N84").Value = -65103
$ws.Range("H113").Value = 734.2632
$ws.Range("J113").Value = 1001.2
$ws.Range("L113").Value = 3003.6
$ws.Range("N113").Value = -7343.6
